$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the API endpoint column (D) for the Type Ahead test rows:
# prefix "/suggest/" -> "/suggest/wos/"
$ws.Range("D3").Value = "/suggest/wos/(S1_TC_T1_hits.hits._source.category)"
$ws.Range("D4").Value = "/suggest/wos/mic*"
$ws.Range("D5").Value = "/suggest/wos/bio methanol"
$ws.Range("D6").Value = "/suggest/wos/bio+methanol"

# Update the view/selection state to match: scroll back to show column A,
# and select cell A6 instead of L2:L6.
$ws.Range("A6").Select()
